# Applies the "Updated cryptos list" data refresh to the cryptos sheet.
# Each assignment below sets a cell to its new post-update value. D-column
# values that look like plain numbers are prefixed with a leading apostrophe
# so Excel keeps storing them as text (matching the original inlineStr cells)
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.634.53"
$ws.Range("E2").Value = "  -2.24%  "

$ws.Range("D3").Value = "3.483.26"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'568.45"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "'182.65"
$ws.Range("E6").Value = "  -3.60%  "

$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -3.31%  "

$ws.Range("D8").Value = "3.474.06"
$ws.Range("E8").Value = "  -3.16%  "

$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").Value = "'0.182"
$ws.Range("E10").Value = "  +1.79%  "

$ws.Range("D11").Value = "'0.636"
$ws.Range("E11").Value = "  -3.85%  "

$ws.Range("D12").Value = "'53.53"
$ws.Range("E12").Value = "  -4.30%  "

$ws.Range("D13").Value = "'0.0000298"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").Value = "'9.39"
$ws.Range("E14").Value = "  -2.61%  "

$ws.Range("D15").Value = "4.046.13"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("D16").Value = "'19.15"
$ws.Range("E16").Value = "  -4.00%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "68.602.65"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.486.22"
$ws.Range("E18").Value = "  -2.82%  "

$ws.Range("D19").Value = "'12.32"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("E20").Value = "  -1.47%  "

$ws.Range("D21").Value = "'538.16"
$ws.Range("E21").Value = "  +13.05%  "

$ws.Range("D22").Value = "'1.01"
$ws.Range("E22").Value = "  -3.09%  "

$ws.Range("D23").Value = "'19.55"
$ws.Range("E23").Value = "  +2.20%  "

$ws.Range("D24").Value = "'4.98"
$ws.Range("E24").Value = "  -2.27%  "

$ws.Range("D25").Value = "'4.38"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").Value = "'93.74"
$ws.Range("E26").Value = "  +3.98%  "

$ws.Range("D27").Value = "'2.90"
$ws.Range("E27").Value = "  -4.64%  "

$ws.Range("D28").Value = "'10.79"
$ws.Range("E28").Value = "  -2.37%  "

$ws.Range("D29").Value = "'9.00"
$ws.Range("E29").Value = "  -3.66%  "

$ws.Range("D30").Value = "'31.28"
$ws.Range("E30").Value = "  -3.01%  "

$ws.Range("D31").Value = "'7.16"
$ws.Range("E31").Value = "  -7.19%  "

$ws.Range("D32").Value = "'12.57"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("D33").Value = "'64.13"
$ws.Range("E33").Value = "  -3.30%  "

$ws.Range("D34").Value = "'0.113"
$ws.Range("E34").Value = "  -5.84%  "

$ws.Range("D35").Value = "'569.63"
$ws.Range("E35").Value = "  -2.79%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'37.71"
$ws.Range("E37").Value = "  -3.58%  "

$ws.Range("D38").Value = "'0.395"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +4.43%  "

$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -4.88%  "

$ws.Range("D41").Value = "'0.132"
$ws.Range("E41").Value = "  -5.40%  "

$ws.Range("D42").Value = "'3.05"
$ws.Range("E42").Value = "  -6.73%  "

$ws.Range("D43").Value = "'3.31"
$ws.Range("E43").Value = "  -5.28%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.191.47"
$ws.Range("E44").Value = "  -1.62%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.95"
$ws.Range("E45").Value = "  -5.07%  "

$ws.Range("D46").Value = "'3.46"
$ws.Range("E46").Value = "  +3.32%  "

$ws.Range("D47").Value = "'0.0435"
$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("D48").Value = "'9.04"
$ws.Range("E48").Value = "  -4.18%  "

$ws.Range("D49").Value = "'0.133"
$ws.Range("E49").Value = "  -3.01%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "'135.75"
$ws.Range("E51").Value = "  -0.93%  "
